$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one
$ws.Range("A1").EntireColumn.Insert()

# Set header and new Lab. # value
$ws.Range("A1").Value = "Lab. #"
$ws.Range("A2").Value = 8166

# Set column A width to match new narrower layout
# (engine quantizes ColumnWidth to the nearest 1/6-character pixel grid; 6.8333
#  lands in the bucket closest to the target stored width of 7.7109375)
$ws.Range("A1").EntireColumn.ColumnWidth = 6.8333333333333335

# Apply the highlight fill color to the full data row (A2:U2)
$ws.Range("A2:U2").Interior.Color = 12379352
